$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Resultado esperado" column (F),
# shifting it (and the blank column after it) one to the right. Insert()
# copies formatting from the split column automatically, so the new F1 and
# the shifted G1 both keep the original header style.
$ws.Columns("F").Insert()

# New column header.
$ws.Range("F1").Value = "Datos"

# Match the original ("Resultado esperado") column's width.
$ws.Columns("F").ColumnWidth = 101.16666666666667

# Re-apply the autofilter over the new, wider data range (toggle off then
# back on so the stored ref actually updates).
$null = $ws.Range("A1:G12").AutoFilter()
$null = $ws.Range("A1:G12").AutoFilter()

# The hidden _FilterDatabase defined name tracks the autofilter range too -
# update it explicitly to the new extent.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Casos de prueba'!`$A`$1:`$G`$12"
    }
}

# Reset the active selection to A1.
$null = $ws.Range("A1").Select()
